# Test_Cases with Hardware and software validation
# Adds a new "Validate Hardware" / "Validate Software" section (TC4.0 / TC4.1)
# to Sheet1, mirroring the existing section layout (blank separator row,
# section-title row, then one or more test-case rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

function Copy-CellFormat($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

function New-BlankSeparatorRow($row) {
    # Mirrors row 17 / row 33: a thin blank divider row.
    Copy-CellFormat "A17:C17" "A$($row):C$($row)"
    Copy-CellFormat "E17:H17" "E$($row):H$($row)"
}

function New-SectionHeaderRow($row, $number, $title) {
    # Mirrors row 18 / row 34: big numbered section title.
    Copy-CellFormat "A34" "A$($row)"
    Copy-CellFormat "C34" "C$($row)"
    $ws.Cells.Item($row, 1).Value = $number
    $ws.Cells.Item($row, 3).Value = $title
}

function New-TestCaseRow($row, $id, $priority, $title, $pre, $steps, $expected, $status) {
    # Mirrors the long wrapped test-case rows (e.g. row 47).
    Copy-CellFormat "A47" "A$($row)"
    Copy-CellFormat "B47" "B$($row)"
    Copy-CellFormat "C47" "C$($row)"
    Copy-CellFormat "D47" "D$($row)"
    Copy-CellFormat "E47" "E$($row)"
    Copy-CellFormat "F47" "F$($row)"
    Copy-CellFormat "G47" "G$($row)"
    Copy-CellFormat "H47" "H$($row)"
    $ws.Rows.Item($row).RowHeight = 409.5

    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 2).Value = $priority
    $ws.Cells.Item($row, 3).Value = $title
    $ws.Cells.Item($row, 4).Value = $pre
    $ws.Cells.Item($row, 5).Value = $steps
    $ws.Cells.Item($row, 6).Value = $expected
    $ws.Cells.Item($row, 8).Value = $status
}

$preConditions = "1- LCD 2*16 is available.`n2- Atmega 32 is available.`n3- 5 volt battary.`n"

$stepsLightAllLeds = "1- Connect pin VSS in LCD to ground.`n2- Connect pin VCC in LCD to 5V battery.`n3- Connect pin VEE in LCD to Resistance 10K then the other leg of the resistance to the ground.`n4- Connect RS pin in LCD to PC0 in ATmega32.`n5- Connect RW pin in LCD to PC1 in ATmega32.`n6- Connect E pin in LCD to PC2 in ATmega32.`n7- Connect PD0 in LCD to pin PB0 in Atmega32.`n8- Connect PD1 in LCD to pin PB1 in Atmega32.`n9- Connect PD2 in LCD to pin PB2 in Atmega32.`n10- Connect PD3 in LCD to pin PB3 in Atmega32.`n11- Connect PD4 in LCD to pin PB4 in Atmega32.`n12- Connect PD5 in LCD to pin PB5 in Atmega32.`n13- Connect PD6 in LCD to pin PB6 in Atmega32.`n14- Connect PD7 in LCD to pin PB7 in Atmega32.`n15- Connect VCC in Atmega32 to 5Volt battarey.`n16- Connect GND pin in Atmega32 to Ground.`n17- make code with software to light all Leds of the LCD."

$stepsMeasureDistanceHw = "1- Connect pin VSS in LCD to ground.`n2- Connect pin VCC in LCD to 5V battery.`n3- Connect pin VEE in LCD to Resistance 10K then the other leg of the resistance to the ground.`n4- Connect RS pin in LCD to PC0 in ATmega32.`n5- Connect RW pin in LCD to PC1 in ATmega32.`n6- Connect E pin in LCD to PC2 in ATmega32.`n7- Connect PD0 in LCD to pin PB0 in Atmega32.`n8- Connect PD1 in LCD to pin PB1 in Atmega32.`n9- Connect PD2 in LCD to pin PB2 in Atmega32.`n10- Connect PD3 in LCD to pin PB3 in Atmega32.`n11- Connect PD4 in LCD to pin PB4 in Atmega32.`n12- Connect PD5 in LCD to pin PB5 in Atmega32.`n13- Connect PD6 in LCD to pin PB6 in Atmega32.`n14- Connect PD7 in LCD to pin PB7 in Atmega32.`n15- Connect VCC in Atmega32 to 5Volt battarey.`n16- Connect GND pin in Atmega32 to Ground.`n17- Connect treg pin in ultrasonic to PA0 in Atmega32.`n18- Connect echo pin in ultrasonic to PD6 in Atmega32.`n19- Connect VCC pin in ultrasonic to 5 volt.`n20- Connect GND pin in ultrasonic to Ground.`n21- make code with software to test the sensor by displaying reading of distance in LCD.`n22- make an object facing the ultrasonic with less than 300cm and more than 2 cm."

$stepsConfigEchoTrigInputOutput = "1- Connect pin VSS in LCD to ground.`n2- Connect pin VCC in LCD to 5V battery.`n3- Connect pin VEE in LCD to Resistance 10K then the other leg of the resistance to the ground.`n4- Connect RS pin in LCD to PC0 in ATmega32.`n5- Connect RW pin in LCD to PC1 in ATmega32.`n6- Connect E pin in LCD to PC2 in ATmega32.`n7- Connect PD0 in LCD to pin PB0 in Atmega32.`n8- Connect PD1 in LCD to pin PB1 in Atmega32.`n9- Connect PD2 in LCD to pin PB2 in Atmega32.`n10- Connect PD3 in LCD to pin PB3 in Atmega32.`n11- Connect PD4 in LCD to pin PB4 in Atmega32.`n12- Connect PD5 in LCD to pin PB5 in Atmega32.`n13- Connect PD6 in LCD to pin PB6 in Atmega32.`n14- Connect PD7 in LCD to pin PB7 in Atmega32.`n15- Connect VCC in Atmega32 to 5Volt battarey.`n16- Connect GND pin in Atmega32 to Ground.`n17- Connect treg pin in ultrasonic to PA0 in Atmega32.`n18- Connect echo pin in ultrasonic to PD6 in Atmega32.`n19- Connect VCC pin in ultrasonic to 5 volt.`n20- Connect GND pin in ultrasonic to Ground.`n21- Configure echo pin in ultrasonic as input pin.`n22- Configure trig pin in ultrasonic as output pin.`n23- make code with software to test the sensor by displaying reading of distance in LCD.`n22- make an object facing the ultrasonic with less than 300cm and more than 2 cm."

$stepsConfigEchoTrigOutputInput = "1- Connect pin VSS in LCD to ground.`n2- Connect pin VCC in LCD to 5V battery.`n3- Connect pin VEE in LCD to Resistance 10K then the other leg of the resistance to the ground.`n4- Connect RS pin in LCD to PC0 in ATmega32.`n5- Connect RW pin in LCD to PC1 in ATmega32.`n6- Connect E pin in LCD to PC2 in ATmega32.`n7- Connect PD0 in LCD to pin PB0 in Atmega32.`n8- Connect PD1 in LCD to pin PB1 in Atmega32.`n9- Connect PD2 in LCD to pin PB2 in Atmega32.`n10- Connect PD3 in LCD to pin PB3 in Atmega32.`n11- Connect PD4 in LCD to pin PB4 in Atmega32.`n12- Connect PD5 in LCD to pin PB5 in Atmega32.`n13- Connect PD6 in LCD to pin PB6 in Atmega32.`n14- Connect PD7 in LCD to pin PB7 in Atmega32.`n15- Connect VCC in Atmega32 to 5Volt battarey.`n16- Connect GND pin in Atmega32 to Ground.`n17- Connect treg pin in ultrasonic to PA0 in Atmega32.`n18- Connect echo pin in ultrasonic to PD6 in Atmega32.`n19- Connect VCC pin in ultrasonic to 5 volt.`n20- Connect GND pin in ultrasonic to Ground.`n21- Configure echo pin in ultrasonic as output pin.`n22- Configure trig pin in ultrasonic as input pin.`n23- make code with software to test the sensor by displaying reading of distance in LCD.`n24- make an object facing the ultrasonic with less than 300cm and more than 2 cm."

# ---- Section 4: Validate Hardware --------------------------------------
New-BlankSeparatorRow 48
New-SectionHeaderRow 49 4 "Validate Hardware "

New-TestCaseRow 50 "TC4.0" "High" `
    "Validate connection of LCD to Atmega32 when make orderusing Atmega to light all leds in LCD" `
    $preConditions $stepsLightAllLeds "All leds in LCD should be lighten." "Open"

New-TestCaseRow 51 "TC4.1" "High" `
    "Validate connection of Ultra sonic to Atmega32 and the LCD when make order using Atmega to Measure the distance " `
    $preConditions $stepsMeasureDistanceHw "the readings should be displayed in LCD " "Open"

# ---- Section 4 (cont'd): Validate Software -------------------------------
New-BlankSeparatorRow 52
New-SectionHeaderRow 53 4 "Validate Software"

New-TestCaseRow 54 "TC4.0" "High" `
    "Validate configuration of echo and trig pins in Ultrasonic sensor" `
    $preConditions $stepsConfigEchoTrigInputOutput "the readings should be displayed in LCD " "Open"

New-TestCaseRow 55 "TC4.0" "High" `
    "Validate configuration of echo and trig pins in Ultrasonic sensor" `
    $preConditions $stepsConfigEchoTrigOutputInput "No reading displayed in LCD." "Open"

# ---- View state: land on the newly added rows ----------------------------
$ws.Application.Goto($ws.Range("D56"), $true)
$ws.Range("D56").Select() | Out-Null
